$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2917716402565462
$ws.Range("C2").Value = 117.745847958593
$ws.Range("D2").Value = 261.3203778131603
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("G2").Value = 2196358.236459397
